$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "pt_min" column (F), shifting F:Q -> G:R
$ws.Columns("F").Insert()

# Header for the newly inserted column F (shared string "eta")
$ws.Range("F1").Value = "eta"

# Populate the new column's data values for rows 2-7
$ws.Range("F2").Value = -1.24
$ws.Range("F3").Value = -0.72
$ws.Range("F4").Value = -0.25
$ws.Range("F5").Value = 0.25
$ws.Range("F6").Value = 0.72
$ws.Range("F7").Value = 1.24

# Match style (centered) used by the rest of the header/data row, same as sibling cells
$ws.Range("F1:F7").HorizontalAlignment = -4108

# Restore the cursor/selection position recorded by Excel after the edit
$ws.Range("F18").Select()

Write-Output "done"
